# New fish-monitoring graphs: add benthic + invert survey rows for the
# "Pillars" (formerly "Pillars of Hercules") and new "Mermaids" sites.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# benthic sheet: append 4 new rows (rows 6-9), copying the date column's
# formatting from the row above so it keeps its date style.
# ---------------------------------------------------------------------
$benthic = $wb.Worksheets.Item("benthic")
$benthic.Activate()

$benthicRows = @(
    @(1,10.1,8.2,28,1,1,1,2),
    @(1,10.1,8.2,28,1,2,1,1),
    @(2,10.1,8.2,28,1,1,3,1),
    @(2,10.1,8.2,28,1,2,2,3)
)

$r = 6
foreach ($row in $benthicRows) {
    $benthic.Cells.Item(5, 3).Copy()
    $benthic.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)

    $benthic.Cells.Item($r, 1).Value = "MWW"
    $benthic.Cells.Item($r, 2).Value = "Pillars of Hercules"
    $benthic.Cells.Item($r, 3).Value = 44477
    $benthic.Cells.Item($r, 4).Value = $row[0]
    $benthic.Cells.Item($r, 5).Value = $row[1]
    $benthic.Cells.Item($r, 6).Value = $row[2]
    $benthic.Cells.Item($r, 7).Value = "m"
    $benthic.Cells.Item($r, 8).Value = $row[3]
    $benthic.Cells.Item($r, 9).Value = "c"
    $benthic.Cells.Item($r, 11).Value = $row[4]
    $benthic.Cells.Item($r, 12).Value = $row[5]
    $benthic.Cells.Item($r, 13).Value = $row[6]
    $benthic.Cells.Item($r, 14).Value = $row[7]
    $r = $r + 1
}

$benthic.Range("M11").Select()

# ---------------------------------------------------------------------
# inverts sheet: add the two new "Diadema" counts for the (still named)
# "Pillars of Hercules" site, then add every "Mermaids" site row, and
# finally rename "Pillars of Hercules" -> "Pillars" everywhere (this
# ordering matches how the new shared strings end up indexed: Diadema,
# Mermaids, Pillars).
# ---------------------------------------------------------------------
$inverts = $wb.Worksheets.Item("inverts")
$inverts.Activate()

$diademaRows = @(
    @(44477,1,3),
    @(44477,2,2)
)
$r = 8
foreach ($row in $diademaRows) {
    $inverts.Cells.Item(7, 3).Copy()
    $inverts.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)

    $inverts.Cells.Item($r, 1).Value = "MWW"
    $inverts.Cells.Item($r, 2).Value = "Pillars of Hercules"
    $inverts.Cells.Item($r, 3).Value = $row[0]
    $inverts.Cells.Item($r, 4).Value = $row[1]
    $inverts.Cells.Item($r, 5).Value = "Diadema"
    $inverts.Cells.Item($r, 6).Value = $row[2]
    $r = $r + 1
}

$mermaidsRows = @(
    @(44842,1,"Lobster",1),
    @(44842,1,"Conch",1),
    @(44842,1,"Lobster",2),
    @(44477,1,"Lobster",1),
    @(44477,1,"Conch",2),
    @(44477,1,"Lobster",3),
    @(44477,1,"Diadema",3),
    @(44477,2,"Diadema",2)
)
foreach ($row in $mermaidsRows) {
    $inverts.Cells.Item(7, 3).Copy()
    $inverts.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)

    $inverts.Cells.Item($r, 1).Value = "MWW"
    $inverts.Cells.Item($r, 2).Value = "Mermaids"
    $inverts.Cells.Item($r, 3).Value = $row[0]
    $inverts.Cells.Item($r, 4).Value = $row[1]
    $inverts.Cells.Item($r, 5).Value = $row[2]
    $inverts.Cells.Item($r, 6).Value = $row[3]
    $r = $r + 1
}

for ($i = 2; $i -le 9; $i++) {
    $inverts.Cells.Item($i, 2).Value = "Pillars"
}
$inverts.Cells.Item(4, 4).Value = 2

$inverts.Range("B15:B17").Select()

# ---------------------------------------------------------------------
# tab / view state: "benthic" is now the active tab
# ---------------------------------------------------------------------
$benthic.Activate()
